$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.045.79"
$ws.Range("E2").Value = "'  +0.14%  "
$ws.Range("D3").Value = "'2.759.13"
$ws.Range("E3").Value = "'  +0.40%  "
$ws.Range("E4").Value = "'  +0.26%  "
$ws.Range("D5").Value = "'579.01"
$ws.Range("E5").Value = "'  +0.35%  "
$ws.Range("D6").Value = "'158.30"
$ws.Range("E6").Value = "'  +3.14%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "'  +0.23%  "
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E9").Value = "'  -1.38%  "
$ws.Range("D10").Value = "'5.78"
$ws.Range("E10").Value = "'  -13.83%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("E11").Value = "'  -1.16%  "
$ws.Range("D12").Value = "'0.158"
$ws.Range("E12").Value = "'  -2.99%  "
$ws.Range("D13").Value = "'3.243.17"
$ws.Range("E13").Value = "'  +0.71%  "
$ws.Range("D14").Value = "'27.06"
$ws.Range("E14").Value = "'  +2.78%  "
$ws.Range("D15").Value = "'63.721.75"
$ws.Range("E15").Value = "'  -0.13%  "
$ws.Range("D16").Value = "'0.0000153"
$ws.Range("E16").Value = "'  +0.23%  "
$ws.Range("D17").Value = "'2.760.90"
$ws.Range("E17").Value = "'  +0.55%  "
$ws.Range("D18").Value = "'12.10"
$ws.Range("E18").Value = "'  +1.27%  "
$ws.Range("D19").Value = "'4.88"
$ws.Range("E19").Value = "'  +0.39%  "
$ws.Range("D20").Value = "'360.07"
$ws.Range("E20").Value = "'  -0.15%  "
$ws.Range("D21").Value = "'6.84"
$ws.Range("E21").Value = "'  -1.41%  "
$ws.Range("D22").Value = "'0.550"
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "'  +0.41%  "
$ws.Range("D24").Value = "'65.78"
$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "'  +0.92%  "
$ws.Range("B26").Value = "'InternetComputer(DFINITY)"
$ws.Range("C26").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'8.49"
$ws.Range("E26").Value = "'  -1.02%  "
$ws.Range("B27").Value = "'Binance-PegBSC-USD"
$ws.Range("C27").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "'  +0.05%  "
$ws.Range("D28").Value = "'0.0₃0931"
$ws.Range("E28").Value = "'  +3.07%  "
$ws.Range("D29").Value = "'1.96"
$ws.Range("E29").Value = "'  -2.32%  "
$ws.Range("D30").Value = "'7.02"
$ws.Range("E30").Value = "'  -1.33%  "
$ws.Range("D31").Value = "'1.24"
$ws.Range("E31").Value = "'  +1.91%  "
$ws.Range("D32").Value = "'167.09"
$ws.Range("E32").Value = "'  -2.86%  "
$ws.Range("D33").Value = "'20.33"
$ws.Range("E33").Value = "'  -0.71%  "
$ws.Range("E34").Value = "'  +3.44%  "
$ws.Range("E35").Value = "'  +0.16%  "
$ws.Range("E36").Value = "'  +1.74%  "
$ws.Range("D37").Value = "'1.80"
$ws.Range("E37").Value = "'  -0.82%  "
$ws.Range("E38").Value = "'  -0.41%  "
$ws.Range("D39").Value = "'6.18"
$ws.Range("E39").Value = "'  +11.47%  "
$ws.Range("D40").Value = "'4.17"
$ws.Range("E40").Value = "'  -1.27%  "
$ws.Range("D41").Value = "'330.21"
$ws.Range("E41").Value = "'  -3.52%  "
$ws.Range("D42").Value = "'39.34"
$ws.Range("E42").Value = "'  +0.61%  "
$ws.Range("E43").Value = "'  -0.91%  "
$ws.Range("B44").Value = "'InjectiveProtocol"
$ws.Range("C44").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'21.85"
$ws.Range("E44").Value = "'  +0.20%  "
$ws.Range("B45").Value = "'Hedera"
$ws.Range("C45").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0594"
$ws.Range("E45").Value = "'  +1.04%  "
$ws.Range("B46").Value = "'VeChain"
$ws.Range("C46").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0257"
$ws.Range("E46").Value = "'  +0.40%  "
$ws.Range("B47").Value = "'Mantle"
$ws.Range("C47").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.634"
$ws.Range("E47").Value = "'  -1.81%  "
$ws.Range("D48").Value = "'135.77"
$ws.Range("E48").Value = "'  -2.50%  "
$ws.Range("E49").Value = "'  +0.32%  "
$ws.Range("E50").Value = "'  -0.15%  "
$ws.Range("D51").Value = "'11.05"
$ws.Range("E51").Value = "'  +0.70%  "
